# Update "Prix Spot" sheet: add a new day column CB (01-sep) with its hourly prices
$wb = $excel.ActiveWorkbook

$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Copy formatting of the last header/date column (CA) into the new one (CB),
# then set the header label.
$wsPrix.Range("CA1").Copy()
$wsPrix.Range("CB1").PasteSpecial(-4122)
$wsPrix.Range("CB1").Value = "01-sep"

# New hourly values for 01-sep
$wsPrix.Range("CB2").Value = 35.79
$wsPrix.Range("CB3").Value = 18.07
$wsPrix.Range("CB4").Value = 22.34
$wsPrix.Range("CB5").Value = 17.64
$wsPrix.Range("CB6").Value = 15
$wsPrix.Range("CB7").Value = 25.21
$wsPrix.Range("CB8").Value = 33.23
$wsPrix.Range("CB9").Value = 53.18
$wsPrix.Range("CB10").Value = 64
$wsPrix.Range("CB11").Value = 35.79
$wsPrix.Range("CB12").Value = 14.03
$wsPrix.Range("CB13").Value = 16.64
$wsPrix.Range("CB14").Value = 13.73
$wsPrix.Range("CB15").Value = 6.07
$wsPrix.Range("CB16").Value = 0
$wsPrix.Range("CB17").Value = 4.65
$wsPrix.Range("CB18").Value = 8.13
$wsPrix.Range("CB19").Value = 14.08
$wsPrix.Range("CB20").Value = 23.08
$wsPrix.Range("CB21").Value = 55
$wsPrix.Range("CB22").Value = 80
$wsPrix.Range("CB23").Value = 84.24
$wsPrix.Range("CB24").Value = 84.74
$wsPrix.Range("CB25").Value = 78.5

# Update "Gaz" sheet: append two new daily rows.
# The "Date" column stores plain text labels (e.g. "2025-08-29"), not real
# dates. Typing an ISO-looking string straight into Value would get it
# auto-recognised as a date, so instead enter it as a text formula
# (="2025-08-30") and immediately collapse it to its static value with a
# values-only paste - this keeps the cell a plain text cell without
# introducing any new cell style/number format.
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A77").Formula = "=""2025-08-30"""
$wsGaz.Range("A77").Copy()
$wsGaz.Range("A77").PasteSpecial(-4163)
$wsGaz.Range("B77").Value = 30.225

$wsGaz.Range("A78").Formula = "=""2025-08-31"""
$wsGaz.Range("A78").Copy()
$wsGaz.Range("A78").PasteSpecial(-4163)
$wsGaz.Range("B78").Value = 30.225

# Update "CO2" sheet: append two new daily rows (same text-date handling).
$wsCo2 = $wb.Worksheets.Item("CO2")

$wsCo2.Range("A77").Formula = "=""2025-08-30"""
$wsCo2.Range("A77").Copy()
$wsCo2.Range("A77").PasteSpecial(-4163)
$wsCo2.Range("B77").Value = 71.1

$wsCo2.Range("A78").Formula = "=""2025-08-31"""
$wsCo2.Range("A78").Copy()
$wsCo2.Range("A78").PasteSpecial(-4163)
$wsCo2.Range("B78").Value = 71.1
